$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999616223267,
    0.99450466596984444,
    0.97302612984653192,
    0.96305391299206267,
    0.95348689418406019,
    0.92996322211706772,
    0.92852461239672435,
    0.92724340695955676,
    0.92898510102145815,
    0.93170499133997331,
    0.93225075419772052,
    0.93354732580627764,
    0.94184701235572743,
    0.9464492878681644,
    0.95207150978729527,
    0.95931112251518358,
    0.96206214793947276,
    0.96394598126022313,
    0.99168302884510617,
    0.96888692145156785,
    0.96271845162810288,
    0.95270081190169442,
    0.96599046650033427,
    0.95296871559541341,
    0.94651156878054854,
    0.92630621526140633,
    0.92145777196945478,
    0.89997870443003536,
    0.88470492511729526,
    0.87813397880391553,
    0.87047965784193582,
    0.86880028122632202,
    0.86828026181625384
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
